$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 111815114
$ws.Range("B16").Value = 90660
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 4362
$ws.Range("F16").Value = "Blå taggsvamp"
$ws.Range("G16").Value = "Hydnellum caeruleum"
$ws.Range("H16").Value = "(Hornem.) P.Karst."
$ws.Range("P16").Value = "åsele 1:1, Ås lm"
$ws.Range("Q16").Value = 610384.0265214761
$ws.Range("R16").Value = 7121170.261031131
$ws.Range("S16").Value = 5
$ws.Range("Z16").Value = "18:19"
$ws.Range("AB16").Value = "18:19"
$ws.Range("A17").Value = 111815269
$ws.Range("B17").Value = 90666
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 4364
$ws.Range("F17").Value = "Dropptaggsvamp"
$ws.Range("G17").Value = "Hydnellum ferrugineum"
$ws.Range("H17").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("P17").Value = "åsele 1:1 (åsele 1:1), Ås lm"
$ws.Range("Q17").Value = 610053.7842541422
$ws.Range("R17").Value = 7121273.15248157
$ws.Range("S17").Value = 1
$ws.Range("Z17").Value = "18:27"
$ws.Range("AB17").Value = "18:27"
$ws.Range("A19").Value = 111814478
$ws.Range("Q19").Value = 610155.3487898401
$ws.Range("R19").Value = 7121461.207019502
$ws.Range("Z19").Value = "17:41"
$ws.Range("AB19").Value = "17:41"
$ws.Range("A20").Value = 111814591
$ws.Range("B20").Value = 77515
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 6425
$ws.Range("F20").Value = "Garnlav"
$ws.Range("G20").Value = "Alectoria sarmentosa"
$ws.Range("H20").Value = "(Ach.) Ach."
$ws.Range("Q20").Value = 610012.4812897337
$ws.Range("R20").Value = 7121464.398116477
$ws.Range("Z20").Value = "17:50"
$ws.Range("AB20").Value = "17:50"
$ws.Range("A21").Value = 111814104
$ws.Range("B21").Value = 56398
$ws.Range("E21").Value = 100109
$ws.Range("F21").Value = "Tretåig hackspett"
$ws.Range("G21").Value = "Picoides tridactylus"
$ws.Range("H21").Value = "(Linnaeus, 1758)"
$ws.Range("Q21").Value = 610154.5078508666
$ws.Range("R21").Value = 7121460.305022033
$ws.Range("Z21").Value = "17:23"
$ws.Range("AB21").Value = "17:23"
$ws.Range("A22").Value = 111814688
$ws.Range("B22").Value = 90087
$ws.Range("D22").Value = "LC"
$ws.Range("E22").Value = 3298
$ws.Range("F22").Value = "Trådticka"
$ws.Range("G22").Value = "Climacocystis borealis"
$ws.Range("H22").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q22").Value = 610011.2059644217
$ws.Range("R22").Value = 7121475.688616944
$ws.Range("Z22").Value = "17:55"
$ws.Range("AB22").Value = "17:55"
$ws.Range("A25").Value = 112013690
$ws.Range("B25").Value = 88489
$ws.Range("E25").Value = 1962
$ws.Range("F25").Value = "Vaddporing"
$ws.Range("G25").Value = "Anomoporia kamtschatica"
$ws.Range("H25").Value = "(Parmasto) Bondartseva"
$ws.Range("Q25").Value = 610051.8565798617
$ws.Range("R25").Value = 7121425.252971379
$ws.Range("Z25").Value = "19:43"
$ws.Range("AB25").Value = "19:43"
$ws.Range("A26").Value = 112013698
$ws.Range("B26").Value = 77515
$ws.Range("E26").Value = 6425
$ws.Range("F26").Value = "Garnlav"
$ws.Range("G26").Value = "Alectoria sarmentosa"
$ws.Range("H26").Value = "(Ach.) Ach."
$ws.Range("Q26").Value = 610094.4326785516
$ws.Range("R26").Value = 7121455.546697079
$ws.Range("Z26").Value = "19:49"
$ws.Range("AB26").Value = "19:49"
$ws.Range("A29").Value = 112013703
$ws.Range("B29").Value = 77515
$ws.Range("E29").Value = 6425
$ws.Range("F29").Value = "Garnlav"
$ws.Range("G29").Value = "Alectoria sarmentosa"
$ws.Range("H29").Value = "(Ach.) Ach."
$ws.Range("Q29").Value = 610144.4332068264
$ws.Range("R29").Value = 7121461.253672058
$ws.Range("Z29").Value = "19:28"
$ws.Range("AB29").Value = "19:28"
$ws.Range("A30").Value = 112013697
$ws.Range("B30").Value = 89423
$ws.Range("E30").Value = 5432
$ws.Range("F30").Value = "Granticka"
$ws.Range("G30").Value = "Porodaedalea chrysoloma"
$ws.Range("H30").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q30").Value = 610102.0736959254
$ws.Range("R30").Value = 7121412.654772604
$ws.Range("A31").Value = 112013699
$ws.Range("Q31").Value = 610068.1736430819
$ws.Range("R31").Value = 7121408.394281525
$ws.Range("Z31").Value = "19:40"
$ws.Range("AB31").Value = "19:40"
$ws.Range("A32").Value = 112013704
$ws.Range("B32").Value = 81248
$ws.Range("E32").Value = 1312
$ws.Range("F32").Value = "Gammelgransskål"
$ws.Range("G32").Value = "Pseudographis pinicola"
$ws.Range("H32").Value = "(Nyl.) Rehm"
$ws.Range("Q32").Value = 610093.591720929
$ws.Range("R32").Value = 7121454.644715369
$ws.Range("A33").Value = 112013700
$ws.Range("B33").Value = 77515
$ws.Range("E33").Value = 6425
$ws.Range("F33").Value = "Garnlav"
$ws.Range("G33").Value = "Alectoria sarmentosa"
$ws.Range("H33").Value = "(Ach.) Ach."
$ws.Range("Q33").Value = 610101.9650201321
$ws.Range("R33").Value = 7121415.702941997
$ws.Range("Z33").Value = "19:35"
$ws.Range("AB33").Value = "19:35"
